# API Error Consolidated List - May 28 release update
# - Adds new USPS IMB / returnReceiptTracking related errors to the
#   "Domestic Labels Outbound" sheet (rows 210-216)
# - Fixes the "Domestic Labels Returns " sheet rows 181-186 so column A
#   correctly reads "Domestic Labels Returns" instead of "Domestic Labels
#   Outbound", and appends the matching new row 187 for that sheet.

$wb  = $excel.ActiveWorkbook
$out = $wb.Worksheets.Item("Domestic Labels Outbound")
$ret = $wb.Worksheets.Item("Domestic Labels Returns ")

# ---------------------------------------------------------------------
# 1) "Domestic Labels Outbound" - new rows 210-216
# ---------------------------------------------------------------------

# Row 210 - full-format copy from an existing similarly formatted row
$out.Range("A68:G68").Copy()
$out.Range("A210:G210").PasteSpecial(-4122)
$out.Rows.Item(210).RowHeight = 30

$out.Range("A210").Value = "Domestic Labels Outbound"
$out.Range("B210").Value = 400
$out.Range("C210").Value = "BAD_REQUEST"
$out.Range("D210").Value = "https://api.usps.com/labels/v3/label"
$out.Range("E210").Value = 160360
$out.Range("F210").Value = "packageOptions.returnReceiptTracking"
$out.Range("G210").Value = "packageOptions.originalPackage.originalTrackingNumber is currently not supported when packageOptions.returnReceiptTracking is true"

# Rows 211-215 - share a common A/B/C/E/F/G style; D carries the new
# indicia URL (unstyled, matching the rest of the plain data cells)
$out.Range("A68:C68").Copy()
$out.Range("A211:C215").PasteSpecial(-4122)
$out.Range("E68").Copy()
$out.Range("E211:E215").PasteSpecial(-4122)
$out.Range("F191").Copy()
$out.Range("F211:F215").PasteSpecial(-4122)
$out.Range("G68").Copy()
$out.Range("G211:G215").PasteSpecial(-4122)

$out.Range("A211:A215").Value = "Domestic Labels Outbound"
$out.Range("B211:B215").Value = 400
$out.Range("C211:C215").Value = "BAD_REQUEST"

$out.Range("D211").Value = "https://api.usps.com/labels/v3/indicia/imb/{imb#}"
$out.Range("D212").Value = "https://api.usps.com/labels/v3/indicia/imb/{imb#}"
$out.Range("D213").Value = "https://api.usps.com/labels/v3/indicia/imb/{imb#}"
$out.Range("D214").Value = "https://api.usps.com/labels/v3/indicia/imb/{imb#}"
$out.Range("D215").Value = "https://api.usps.com/labels/v3/indicia/imb/{imb#}"

$out.Range("E211").Value = 160361
$out.Range("E212").Value = 160362
$out.Range("E213").Value = 160363
$out.Range("E214").Value = 160364
$out.Range("E215").Value = 160364

$out.Range("F211").Value = "imb"
$out.Range("G211").Value = "Indicia not found."

$out.Range("F212").Value = "imb.mailingDate"
$out.Range("G212").Value = "letter/flat cannot be more than %d days in the past"

$out.Range("F213").Value = "imb.status"
$out.Range("G213").Value = "letter/flat is already canceled"

$out.Range("F214").Value = "imb.status"
$out.Range("G214").Value = "letter/flat is already disputed"

$out.Range("F215").Value = "imb"
$out.Range("G215").Value = "Forbidden"

# Row 216 - back to the 'imageInfo.labelType' style used by row 210,
# but the G cell gets the new small plain-left-aligned style
$out.Range("A210:F210").Copy()
$out.Range("A216:F216").PasteSpecial(-4122)
$out.Rows.Item(216).RowHeight = 15

$out.Range("A216").Value = "Domestic Labels Outbound"
$out.Range("B216").Value = 400
$out.Range("C216").Value = "BAD_REQUEST"
$out.Range("D216").Value = "https://api.usps.com/labels/v3/label"
$out.Range("E216").Value = 160366
$out.Range("F216").Value = "imageInfo.labelType"

$out.Range("G216").Value = "4X4LABEL currently does not support imageInfo.imageType of 'ZPL203DPI' or 'ZPL300DPI'"
$out.Range("G216").Font.Name = "Calibri"
$out.Range("G216").Font.Size = 11
$out.Range("G216").HorizontalAlignment = -4131

# ---------------------------------------------------------------------
# 2) "Domestic Labels Returns " - fix rows 181-186 (wrong sheet label)
# ---------------------------------------------------------------------
$ret.Range("A180").Copy()
$ret.Range("A181:A186").PasteSpecial(-4122)
$ret.Range("A181:A186").Value = "Domestic Labels Returns"

# ---------------------------------------------------------------------
# 3) "Domestic Labels Returns " - new row 187 (mirrors outbound row 216)
# ---------------------------------------------------------------------
$ret.Range("A186:E186").Copy()
$ret.Range("A187:E187").PasteSpecial(-4122)

$out.Range("F191").Copy()
$ret.Range("F187").PasteSpecial(-4122)

$out.Range("G216").Copy()
$ret.Range("G187").PasteSpecial(-4122)

$ret.Range("A187").Value = "Domestic Labels Returns"
$ret.Range("B187").Value = 400
$ret.Range("C187").Value = "BAD_REQUEST"
$ret.Range("D187").Value = "https://api.usps.com/labels/v3/return-label"
$ret.Range("E187").Value = 160366
$ret.Range("F187").Value = "imageInfo.labelType"
$ret.Range("G187").Value = "4X4LABEL currently does not support imageInfo.imageType of 'ZPL203DPI' or 'ZPL300DPI'"
